$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in two "Razon social" entries (commas -> periods, drop dots in "S.H.")
$ws.Cells.Item(204, 5).Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Cells.Item(215, 5).Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# Normalize "Importe" text values from es-AR grouping (1.234,56) to plain decimal (1234.56).
# Mark the range as Text first so Excel keeps storing these as strings instead of
# silently parsing them back into numbers, then clear the temporary format so the
# cells keep their original (default) style.
$importeRange = $ws.Range("H2:H285")
$importeRange.NumberFormat = "@"

$ws.Cells.Item(2, 8).Value = "27800.00"
$ws.Cells.Item(3, 8).Value = "10500.00"
$ws.Cells.Item(4, 8).Value = "1567200.00"
$ws.Cells.Item(5, 8).Value = "101169.00"
$ws.Cells.Item(6, 8).Value = "433000.00"
$ws.Cells.Item(7, 8).Value = "1204000.00"
$ws.Cells.Item(8, 8).Value = "70.00"
$ws.Cells.Item(9, 8).Value = "11900.00"
$ws.Cells.Item(10, 8).Value = "4500.00"
$ws.Cells.Item(11, 8).Value = "130.00"
$ws.Cells.Item(12, 8).Value = "2315.00"
$ws.Cells.Item(13, 8).Value = "12946.00"
$ws.Cells.Item(14, 8).Value = "1553.20"
$ws.Cells.Item(15, 8).Value = "1000.00"
$ws.Cells.Item(16, 8).Value = "3050.00"
$ws.Cells.Item(17, 8).Value = "952000.00"
$ws.Cells.Item(18, 8).Value = "2900000.00"
$ws.Cells.Item(19, 8).Value = "3398.00"
$ws.Cells.Item(20, 8).Value = "2159.55"
$ws.Cells.Item(21, 8).Value = "3920.00"
$ws.Cells.Item(22, 8).Value = "5379.02"
$ws.Cells.Item(23, 8).Value = "615922.50"
$ws.Cells.Item(24, 8).Value = "533274.30"
$ws.Cells.Item(25, 8).Value = "113360.00"
$ws.Cells.Item(26, 8).Value = "3700.00"
$ws.Cells.Item(27, 8).Value = "124543.30"
$ws.Cells.Item(28, 8).Value = "679781.75"
$ws.Cells.Item(29, 8).Value = "192982.68"
$ws.Cells.Item(30, 8).Value = "3059.85"
$ws.Cells.Item(31, 8).Value = "49530.00"
$ws.Cells.Item(32, 8).Value = "14000.00"
$ws.Cells.Item(33, 8).Value = "1050.00"
$ws.Cells.Item(34, 8).Value = "1931.10"
$ws.Cells.Item(35, 8).Value = "80939.12"
$ws.Cells.Item(36, 8).Value = "17017.00"
$ws.Cells.Item(37, 8).Value = "5230.65"
$ws.Cells.Item(38, 8).Value = "71985.00"
$ws.Cells.Item(39, 8).Value = "70500.00"
$ws.Cells.Item(40, 8).Value = "12800.00"
$ws.Cells.Item(41, 8).Value = "1400.00"
$ws.Cells.Item(42, 8).Value = "19500.00"
$ws.Cells.Item(43, 8).Value = "19123.50"
$ws.Cells.Item(44, 8).Value = "433.49"
$ws.Cells.Item(45, 8).Value = "1780.00"
$ws.Cells.Item(46, 8).Value = "5213.69"
$ws.Cells.Item(47, 8).Value = "3018220.15"
$ws.Cells.Item(48, 8).Value = "17000.81"
$ws.Cells.Item(49, 8).Value = "6960.00"
$ws.Cells.Item(50, 8).Value = "1698770.80"
$ws.Cells.Item(51, 8).Value = "10000.00"
$ws.Cells.Item(52, 8).Value = "17.65"
$ws.Cells.Item(53, 8).Value = "1500.00"
$ws.Cells.Item(54, 8).Value = "4096.70"
$ws.Cells.Item(55, 8).Value = "15817.00"
$ws.Cells.Item(56, 8).Value = "90960.00"
$ws.Cells.Item(57, 8).Value = "14600.00"
$ws.Cells.Item(58, 8).Value = "97939.25"
$ws.Cells.Item(59, 8).Value = "19788.09"
$ws.Cells.Item(60, 8).Value = "2250.00"
$ws.Cells.Item(61, 8).Value = "19148.00"
$ws.Cells.Item(62, 8).Value = "7600.00"
$ws.Cells.Item(63, 8).Value = "88067.28"
$ws.Cells.Item(64, 8).Value = "3740.00"
$ws.Cells.Item(65, 8).Value = "1090.00"
$ws.Cells.Item(66, 8).Value = "5900.00"
$ws.Cells.Item(67, 8).Value = "27118.23"
$ws.Cells.Item(68, 8).Value = "3600.00"
$ws.Cells.Item(69, 8).Value = "5824.90"
$ws.Cells.Item(70, 8).Value = "25737.00"
$ws.Cells.Item(71, 8).Value = "825.00"
$ws.Cells.Item(72, 8).Value = "7596.05"
$ws.Cells.Item(73, 8).Value = "60.00"
$ws.Cells.Item(74, 8).Value = "791093.00"
$ws.Cells.Item(75, 8).Value = "1000.00"
$ws.Cells.Item(76, 8).Value = "9848.00"
$ws.Cells.Item(77, 8).Value = "14850.00"
$ws.Cells.Item(78, 8).Value = "32929.49"
$ws.Cells.Item(79, 8).Value = "3049.00"
$ws.Cells.Item(80, 8).Value = "142000.00"
$ws.Cells.Item(81, 8).Value = "695.87"
$ws.Cells.Item(82, 8).Value = "10811.38"
$ws.Cells.Item(83, 8).Value = "5340.24"
$ws.Cells.Item(84, 8).Value = "2224.80"
$ws.Cells.Item(85, 8).Value = "2500.00"
$ws.Cells.Item(86, 8).Value = "25580.00"
$ws.Cells.Item(87, 8).Value = "17700.00"
$ws.Cells.Item(88, 8).Value = "334200.00"
$ws.Cells.Item(89, 8).Value = "970.00"
$ws.Cells.Item(90, 8).Value = "4785.00"
$ws.Cells.Item(91, 8).Value = "7005.00"
$ws.Cells.Item(92, 8).Value = "39635.01"
$ws.Cells.Item(93, 8).Value = "7500.00"
$ws.Cells.Item(94, 8).Value = "5400.00"
$ws.Cells.Item(95, 8).Value = "231000.00"
$ws.Cells.Item(96, 8).Value = "11400.00"
$ws.Cells.Item(97, 8).Value = "900.00"
$ws.Cells.Item(98, 8).Value = "9900.00"
$ws.Cells.Item(99, 8).Value = "6200.00"
$ws.Cells.Item(100, 8).Value = "18070.00"
$ws.Cells.Item(101, 8).Value = "20160.00"
$ws.Cells.Item(102, 8).Value = "564.48"
$ws.Cells.Item(103, 8).Value = "7500.00"
$ws.Cells.Item(104, 8).Value = "16200.00"
$ws.Cells.Item(105, 8).Value = "9840.00"
$ws.Cells.Item(106, 8).Value = "60.12"
$ws.Cells.Item(107, 8).Value = "5.00"
$ws.Cells.Item(108, 8).Value = "16238.78"
$ws.Cells.Item(109, 8).Value = "18410.00"
$ws.Cells.Item(110, 8).Value = "19470.06"
$ws.Cells.Item(111, 8).Value = "7900.00"
$ws.Cells.Item(112, 8).Value = "3636.79"
$ws.Cells.Item(113, 8).Value = "85419.94"
$ws.Cells.Item(114, 8).Value = "430.00"
$ws.Cells.Item(115, 8).Value = "352590.00"
$ws.Cells.Item(116, 8).Value = "1200.00"
$ws.Cells.Item(117, 8).Value = "284019.01"
$ws.Cells.Item(118, 8).Value = "1710.00"
$ws.Cells.Item(119, 8).Value = "1200.00"
$ws.Cells.Item(120, 8).Value = "33926.38"
$ws.Cells.Item(121, 8).Value = "25631.80"
$ws.Cells.Item(122, 8).Value = "151600.00"
$ws.Cells.Item(123, 8).Value = "93250.00"
$ws.Cells.Item(124, 8).Value = "12850.00"
$ws.Cells.Item(125, 8).Value = "4747.50"
$ws.Cells.Item(126, 8).Value = "22980.00"
$ws.Cells.Item(127, 8).Value = "2646.12"
$ws.Cells.Item(128, 8).Value = "8184.00"
$ws.Cells.Item(129, 8).Value = "66323.23"
$ws.Cells.Item(130, 8).Value = "135.00"
$ws.Cells.Item(131, 8).Value = "320.00"
$ws.Cells.Item(132, 8).Value = "22596.50"
$ws.Cells.Item(133, 8).Value = "15000.00"
$ws.Cells.Item(134, 8).Value = "17760.00"
$ws.Cells.Item(135, 8).Value = "2300.00"
$ws.Cells.Item(136, 8).Value = "3653.55"
$ws.Cells.Item(137, 8).Value = "5452.00"
$ws.Cells.Item(138, 8).Value = "17400.00"
$ws.Cells.Item(139, 8).Value = "2021.54"
$ws.Cells.Item(140, 8).Value = "39200.00"
$ws.Cells.Item(141, 8).Value = "7200.00"
$ws.Cells.Item(142, 8).Value = "6800.00"
$ws.Cells.Item(143, 8).Value = "7311.79"
$ws.Cells.Item(144, 8).Value = "5498.07"
$ws.Cells.Item(145, 8).Value = "15900.00"
$ws.Cells.Item(146, 8).Value = "197087.00"
$ws.Cells.Item(147, 8).Value = "3337.04"
$ws.Cells.Item(148, 8).Value = "2500.00"
$ws.Cells.Item(149, 8).Value = "22000.00"
$ws.Cells.Item(150, 8).Value = "309400.00"
$ws.Cells.Item(151, 8).Value = "1500.00"
$ws.Cells.Item(152, 8).Value = "19000.00"
$ws.Cells.Item(153, 8).Value = "17000.00"
$ws.Cells.Item(154, 8).Value = "22000.00"
$ws.Cells.Item(155, 8).Value = "72500.00"
$ws.Cells.Item(156, 8).Value = "90000.00"
$ws.Cells.Item(157, 8).Value = "3500.00"
$ws.Cells.Item(158, 8).Value = "285900.00"
$ws.Cells.Item(159, 8).Value = "1219510.00"
$ws.Cells.Item(160, 8).Value = "70000.00"
$ws.Cells.Item(161, 8).Value = "153600.00"
$ws.Cells.Item(162, 8).Value = "1767.99"
$ws.Cells.Item(163, 8).Value = "1440.00"
$ws.Cells.Item(164, 8).Value = "37.00"
$ws.Cells.Item(165, 8).Value = "26540.00"
$ws.Cells.Item(166, 8).Value = "58.00"
$ws.Cells.Item(167, 8).Value = "7686.46"
$ws.Cells.Item(168, 8).Value = "35000.00"
$ws.Cells.Item(169, 8).Value = "17700.00"
$ws.Cells.Item(170, 8).Value = "22000.00"
$ws.Cells.Item(171, 8).Value = "26000.00"
$ws.Cells.Item(172, 8).Value = "2200.00"
$ws.Cells.Item(173, 8).Value = "14000.00"
$ws.Cells.Item(174, 8).Value = "12000.00"
$ws.Cells.Item(175, 8).Value = "18000.00"
$ws.Cells.Item(176, 8).Value = "2000.00"
$ws.Cells.Item(177, 8).Value = "30000.00"
$ws.Cells.Item(178, 8).Value = "13000.00"
$ws.Cells.Item(179, 8).Value = "12400.00"
$ws.Cells.Item(180, 8).Value = "9000.00"
$ws.Cells.Item(181, 8).Value = "10200.00"
$ws.Cells.Item(182, 8).Value = "4000.00"
$ws.Cells.Item(183, 8).Value = "14000.00"
$ws.Cells.Item(184, 8).Value = "2500.00"
$ws.Cells.Item(185, 8).Value = "33200.00"
$ws.Cells.Item(186, 8).Value = "10000.00"
$ws.Cells.Item(187, 8).Value = "10000.00"
$ws.Cells.Item(188, 8).Value = "5000.00"
$ws.Cells.Item(189, 8).Value = "73206.25"
$ws.Cells.Item(190, 8).Value = "18000.00"
$ws.Cells.Item(191, 8).Value = "10000.00"
$ws.Cells.Item(192, 8).Value = "16000.00"
$ws.Cells.Item(193, 8).Value = "13000.00"
$ws.Cells.Item(194, 8).Value = "297950.00"
$ws.Cells.Item(195, 8).Value = "18000.00"
$ws.Cells.Item(196, 8).Value = "20000.00"
$ws.Cells.Item(197, 8).Value = "10000.00"
$ws.Cells.Item(198, 8).Value = "10000.00"
$ws.Cells.Item(199, 8).Value = "111300.00"
$ws.Cells.Item(200, 8).Value = "35000.00"
$ws.Cells.Item(201, 8).Value = "11800.00"
$ws.Cells.Item(202, 8).Value = "19600.00"
$ws.Cells.Item(203, 8).Value = "617.43"
$ws.Cells.Item(204, 8).Value = "8660.00"
$ws.Cells.Item(205, 8).Value = "25064.00"
$ws.Cells.Item(206, 8).Value = "36350.00"
$ws.Cells.Item(207, 8).Value = "22300.00"
$ws.Cells.Item(208, 8).Value = "13128.60"
$ws.Cells.Item(209, 8).Value = "200.00"
$ws.Cells.Item(210, 8).Value = "7150.00"
$ws.Cells.Item(211, 8).Value = "10929.87"
$ws.Cells.Item(212, 8).Value = "4163.00"
$ws.Cells.Item(213, 8).Value = "7820.05"
$ws.Cells.Item(214, 8).Value = "29627.60"
$ws.Cells.Item(215, 8).Value = "10780.00"
$ws.Cells.Item(216, 8).Value = "838.40"
$ws.Cells.Item(217, 8).Value = "7893.00"
$ws.Cells.Item(218, 8).Value = "3000.00"
$ws.Cells.Item(219, 8).Value = "4440.00"
$ws.Cells.Item(220, 8).Value = "487.48"
$ws.Cells.Item(221, 8).Value = "1350.00"
$ws.Cells.Item(222, 8).Value = "2585.22"
$ws.Cells.Item(223, 8).Value = "5554.27"
$ws.Cells.Item(224, 8).Value = "5750.00"
$ws.Cells.Item(225, 8).Value = "5817.80"
$ws.Cells.Item(226, 8).Value = "62231.14"
$ws.Cells.Item(227, 8).Value = "12464.09"
$ws.Cells.Item(228, 8).Value = "80000.00"
$ws.Cells.Item(229, 8).Value = "40000.00"
$ws.Cells.Item(230, 8).Value = "40000.00"
$ws.Cells.Item(231, 8).Value = "40000.00"
$ws.Cells.Item(232, 8).Value = "80000.00"
$ws.Cells.Item(233, 8).Value = "40000.00"
$ws.Cells.Item(234, 8).Value = "55000.00"
$ws.Cells.Item(235, 8).Value = "40000.00"
$ws.Cells.Item(236, 8).Value = "40000.00"
$ws.Cells.Item(237, 8).Value = "80000.00"
$ws.Cells.Item(238, 8).Value = "80000.00"
$ws.Cells.Item(239, 8).Value = "151500.00"
$ws.Cells.Item(240, 8).Value = "183348.00"
$ws.Cells.Item(241, 8).Value = "11600.00"
$ws.Cells.Item(242, 8).Value = "7753543.08"
$ws.Cells.Item(243, 8).Value = "30350.00"
$ws.Cells.Item(244, 8).Value = "11000.00"
$ws.Cells.Item(245, 8).Value = "15443949.00"
$ws.Cells.Item(246, 8).Value = "323950.00"
$ws.Cells.Item(247, 8).Value = "325710.00"
$ws.Cells.Item(248, 8).Value = "298100.00"
$ws.Cells.Item(249, 8).Value = "317700.00"
$ws.Cells.Item(250, 8).Value = "298100.00"
$ws.Cells.Item(251, 8).Value = "322600.00"
$ws.Cells.Item(252, 8).Value = "566600.00"
$ws.Cells.Item(253, 8).Value = "298100.00"
$ws.Cells.Item(254, 8).Value = "732700.00"
$ws.Cells.Item(255, 8).Value = "737000.00"
$ws.Cells.Item(256, 8).Value = "746940.00"
$ws.Cells.Item(257, 8).Value = "298100.00"
$ws.Cells.Item(258, 8).Value = "298100.00"
$ws.Cells.Item(259, 8).Value = "596200.00"
$ws.Cells.Item(260, 8).Value = "531300.00"
$ws.Cells.Item(261, 8).Value = "611400.00"
$ws.Cells.Item(262, 8).Value = "885500.00"
$ws.Cells.Item(263, 8).Value = "566600.00"
$ws.Cells.Item(264, 8).Value = "897140.00"
$ws.Cells.Item(265, 8).Value = "610480.00"
$ws.Cells.Item(266, 8).Value = "315520.00"
$ws.Cells.Item(267, 8).Value = "100000.00"
$ws.Cells.Item(268, 8).Value = "4329000.00"
$ws.Cells.Item(269, 8).Value = "19440.00"
$ws.Cells.Item(270, 8).Value = "488000.00"
$ws.Cells.Item(271, 8).Value = "288950.00"
$ws.Cells.Item(272, 8).Value = "40805.00"
$ws.Cells.Item(273, 8).Value = "148800.00"
$ws.Cells.Item(274, 8).Value = "280000.00"
$ws.Cells.Item(275, 8).Value = "14000.00"
$ws.Cells.Item(276, 8).Value = "80027.00"
$ws.Cells.Item(277, 8).Value = "2500.00"
$ws.Cells.Item(278, 8).Value = "2813.38"
$ws.Cells.Item(279, 8).Value = "3300.00"
$ws.Cells.Item(280, 8).Value = "73900.00"
$ws.Cells.Item(281, 8).Value = "13680.00"
$ws.Cells.Item(282, 8).Value = "7000.00"
$ws.Cells.Item(283, 8).Value = "3445.00"
$ws.Cells.Item(284, 8).Value = "11400.00"
$ws.Cells.Item(285, 8).Value = "20200.00"

$importeRange.ClearFormats()

